$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the touched cells keep their original text formatting (values such as
# "0.570" or "2.028.31" are not valid numbers and must stay literal text).
$cells = @('D2', 'E2', 'D3', 'E3', 'E4', 'D5', 'E5', 'D6', 'E6', 'E7', 'D8', 'E8', 'D9', 'E9', 'E10', 'E11', 'E12', 'D13', 'E13', 'D14', 'E14', 'B15', 'C15', 'D15', 'E15', 'B16', 'C16', 'D16', 'E16', 'B17', 'C17', 'D17', 'E17', 'D18', 'E18', 'D19', 'E19', 'E20', 'D21', 'E21', 'E22', 'D23', 'E23', 'E24', 'E25', 'D26', 'E27', 'E28', 'E29', 'E31', 'E32', 'D33', 'E33', 'D34', 'E34', 'E35', 'E36', 'D37', 'E37', 'D38', 'E38', 'D39', 'E39', 'E40', 'E41', 'D42', 'E42', 'E43', 'D44', 'E44', 'B45', 'C45', 'D45', 'E45', 'B46', 'C46', 'D46', 'E46', 'D47', 'E47', 'D48', 'E48', 'D49', 'E49', 'D50', 'E50', 'E51')
foreach ($cellRef in $cells) { $ws.Range($cellRef).NumberFormat = "@" }

$ws.Range('D2').Value = '59.269.29'
$ws.Range('E2').Value = '  +0.31%  '
$ws.Range('D3').Value = '2.636.74'
$ws.Range('E3').Value = '  -0.69%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '530.21'
$ws.Range('E5').Value = '  +1.30%  '
$ws.Range('D6').Value = '145.06'
$ws.Range('E6').Value = '  +0.07%  '
$ws.Range('E7').Value = '  -0.14%  '
$ws.Range('D8').Value = '0.570'
$ws.Range('E8').Value = '  -0.40%  '
$ws.Range('D9').Value = '6.66'
$ws.Range('E9').Value = '  -5.37%  '
$ws.Range('E10').Value = '  +0.34%  '
$ws.Range('E11').Value = '  -0.23%  '
$ws.Range('E12').Value = '  +0.57%  '
$ws.Range('D13').Value = '3.103.90'
$ws.Range('E13').Value = '  -0.59%  '
$ws.Range('D14').Value = '59.251.13'
$ws.Range('E14').Value = '  +0.25%  '
$ws.Range('B15').Value = 'Avalanche'
$ws.Range('C15').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D15').Value = '20.61'
$ws.Range('E15').Value = '  -2.95%  '
$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').Value = '0.0000135'
$ws.Range('E16').Value = '  -0.85%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '2.617.99'
$ws.Range('E17').Value = '  -1.79%  '
$ws.Range('D18').Value = '343.08'
$ws.Range('E18').Value = '  +0.76%  '
$ws.Range('D19').Value = '4.41'
$ws.Range('E19').Value = '  +0.43%  '
$ws.Range('E20').Value = '  +1.13%  '
$ws.Range('D21').Value = '6.33'
$ws.Range('E21').Value = '  -0.60%  '
$ws.Range('E22').Value = '  +0.28%  '
$ws.Range('D23').Value = '66.55'
$ws.Range('E23').Value = '  +4.03%  '
$ws.Range('E24').Value = '  -0.39%  '
$ws.Range('E25').Value = '  +0.08%  '
$ws.Range('D26').Value = '2.757.71'
$ws.Range('E27').Value = '  -0.19%  '
$ws.Range('E28').Value = '  +0.09%  '
$ws.Range('E29').Value = '  -2.06%  '
$ws.Range('E31').Value = '  -6.02%  '
$ws.Range('E32').Value = '  +0.81%  '
$ws.Range('D33').Value = '18.95'
$ws.Range('E33').Value = '  +0.64%  '
$ws.Range('D34').Value = '150.02'
$ws.Range('E34').Value = '  +0.65%  '
$ws.Range('E35').Value = '  -1.24%  '
$ws.Range('E36').Value = '  -3.14%  '
$ws.Range('D37').Value = '36.36'
$ws.Range('E37').Value = '  -0.85%  '
$ws.Range('D38').Value = '0.836'
$ws.Range('E38').Value = '  -7.31%  '
$ws.Range('D39').Value = '0.838'
$ws.Range('E39').Value = '  -5.15%  '
$ws.Range('E40').Value = '  -2.41%  '
$ws.Range('E41').Value = '  -0.16%  '
$ws.Range('D42').Value = '0.998'
$ws.Range('E42').Value = '  -0.10%  '
$ws.Range('E43').Value = '  +0.09%  '
$ws.Range('D44').Value = '0.597'
$ws.Range('E44').Value = '  -3.54%  '
$ws.Range('B45').Value = 'Bittensor'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D45').Value = '268.84'
$ws.Range('E45').Value = '  -2.43%  '
$ws.Range('B46').Value = 'WhiteBITCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D46').Value = '10.72'
$ws.Range('E46').Value = '  +1.75%  '
$ws.Range('D47').Value = '19.05'
$ws.Range('E47').Value = '  -4.20%  '
$ws.Range('D48').Value = '0.0531'
$ws.Range('E48').Value = '  -0.81%  '
$ws.Range('D49').Value = '2.028.31'
$ws.Range('E49').Value = '  -0.06%  '
$ws.Range('D50').Value = '4.69'
$ws.Range('E50').Value = '  -1.75%  '
$ws.Range('E51').Value = '  -0.57%  '
